# Student Management - Academic Year Code Implementation
#
# Update the single permission-slip record on the active sheet:
#   - ParentName (G2):          Robert K. Doe   -> Michael B. Doe
#   - RequestedDateTime (H2):   23-02-2025 ...  -> 26-03-2025 at 08:13 PM
#   - Reason (I2):              Rimjhim is at home -> Student requires early
#                                pickup due to health concerns
#   - PickedUp (J2):            Rahul -> Meera Kapoor
#
# The columns that hold the longer replacement text (ParentName, Reason,
# PickedUp) are then widened so the new values are fully visible, matching
# the wider column layout shipped with this change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Michael B. Doe"
$ws.Range("H2").Value = "26-03-2025 at 08:13 PM"
$ws.Range("I2").Value = "Student requires early pickup due to health concerns"
$ws.Range("J2").Value = "Meera Kapoor"

# Widen ParentName (G), Reason (I) and PickedUp (J) columns to fit the new,
# longer content (RequestedDateTime / column H keeps its existing width).
$ws.Columns.Item(7).ColumnWidth = 13.083333333333334
$ws.Columns.Item(9).ColumnWidth = 46.25
$ws.Columns.Item(10).ColumnWidth = 12.416666666666666
